$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill team record values (Wins=74, Losses=88, Ties=0) for every data row (2-50)
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 74   # AD
    $ws.Cells.Item($r, 31).Value = 88   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
